$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First Rating")

# Copy formatting from the last existing row (277) down through the new rows (278-302)
$ws.Range("A277:G277").Copy()
$ws.Range("A278:G302").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A278").Value = "Château Pavie"
$ws.Range("B278").Value = "Other"
$ws.Range("C278").Value = 1994
$ws.Range("D278").Value = 85
$ws.Range("E278").Value = 87
$ws.Range("F278").Value2 = 34790
$ws.Range("G278").Value = "Robert Parker"

$ws.Range("A279").Value = "Château Pavie"
$ws.Range("B279").Value = "Other"
$ws.Range("C279").Value = 1995
$ws.Range("D279").Value = 85
$ws.Range("E279").Value = 87
$ws.Range("F279").Value2 = 35156
$ws.Range("G279").Value = "Robert Parker"

$ws.Range("A280").Value = "Château Pavie"
$ws.Range("B280").Value = "Other"
$ws.Range("C280").Value = 1996
$ws.Range("D280").Value = 85
$ws.Range("E280").Value = 87
$ws.Range("F280").Value2 = 35796
$ws.Range("G280").Value = "Robert Parker"

$ws.Range("A281").Value = "Château Pavie"
$ws.Range("B281").Value = "Other"
$ws.Range("C281").Value = 1997
$ws.Range("D281").Value = 83
$ws.Range("E281").Value = 85
$ws.Range("F281").Value2 = 35886
$ws.Range("G281").Value = "Robert Parker"

$ws.Range("A282").Value = "Château Pavie"
$ws.Range("B282").Value = "Other"
$ws.Range("C282").Value = 1998
$ws.Range("D282").Value = 91
$ws.Range("E282").Value = 93
$ws.Range("F282").Value2 = 36251
$ws.Range("G282").Value = "Robert Parker"

$ws.Range("A283").Value = "Château Pavie"
$ws.Range("B283").Value = "Other"
$ws.Range("C283").Value = 1999
$ws.Range("D283").Value = 92
$ws.Range("E283").Value = 96
$ws.Range("F283").Value2 = 36617
$ws.Range("G283").Value = "Robert Parker"

$ws.Range("A284").Value = "Château Pavie"
$ws.Range("B284").Value = "Other"
$ws.Range("C284").Value = 2000
$ws.Range("D284").Value = 95
$ws.Range("E284").Value = 96
$ws.Range("F284").Value2 = 36982
$ws.Range("G284").Value = "Robert Parker"

$ws.Range("A285").Value = "Château Pavie"
$ws.Range("B285").Value = "Other"
$ws.Range("C285").Value = 2001
$ws.Range("D285").Value = 94
$ws.Range("E285").Value = 96
$ws.Range("F285").Value2 = 37347
$ws.Range("G285").Value = "Robert Parker"

$ws.Range("A286").Value = "Château Pavie"
$ws.Range("B286").Value = "Other"
$ws.Range("C286").Value = 2002
$ws.Range("D286").Value = 92
$ws.Range("E286").Value = 95
$ws.Range("F286").Value2 = 37895
$ws.Range("G286").Value = "Robert Parker"

$ws.Range("A287").Value = "Château Pavie"
$ws.Range("B287").Value = "Other"
$ws.Range("C287").Value = 2003
$ws.Range("D287").Value = 96
$ws.Range("E287").Value = 100
$ws.Range("F287").Value2 = 38078
$ws.Range("G287").Value = "Robert Parker"

$ws.Range("A288").Value = "Château Pavie"
$ws.Range("B288").Value = "Other"
$ws.Range("C288").Value = 2004
$ws.Range("D288").Value = 95
$ws.Range("E288").Value = 97
$ws.Range("F288").Value2 = 38443
$ws.Range("G288").Value = "Robert Parker"

$ws.Range("A289").Value = "Château Pavie"
$ws.Range("B289").Value = "Other"
$ws.Range("C289").Value = 2005
$ws.Range("D289").Value = 98
$ws.Range("E289").Value = 100
$ws.Range("F289").Value2 = 38808
$ws.Range("G289").Value = "Robert Parker"

$ws.Range("A290").Value = "Château Pavie"
$ws.Range("B290").Value = "Other"
$ws.Range("C290").Value = 2006
$ws.Range("D290").Value = 96
$ws.Range("E290").Value = 98
$ws.Range("F290").Value2 = 39173
$ws.Range("G290").Value = "Robert Parker"

$ws.Range("A291").Value = "Château Pavie"
$ws.Range("B291").Value = "Other"
$ws.Range("C291").Value = 2007
$ws.Range("D291").Value = 93
$ws.Range("E291").Value = 95
$ws.Range("F291").Value2 = 39539
$ws.Range("G291").Value = "Robert Parker"

$ws.Range("A292").Value = "Château Pavie"
$ws.Range("B292").Value = "Other"
$ws.Range("C292").Value = 2008
$ws.Range("D292").Value = 96
$ws.Range("E292").Value = 98
$ws.Range("F292").Value2 = 39904
$ws.Range("G292").Value = "Robert Parker"

$ws.Range("A293").Value = "Château Pavie"
$ws.Range("B293").Value = "Other"
$ws.Range("C293").Value = 2009
$ws.Range("D293").Value = 96
$ws.Range("E293").Value = 100
$ws.Range("F293").Value2 = 40269
$ws.Range("G293").Value = "Robert Parker"

$ws.Range("A294").Value = "Château Pavie"
$ws.Range("B294").Value = "Other"
$ws.Range("C294").Value = 2010
$ws.Range("D294").Value = 95
$ws.Range("E294").Value = 98
$ws.Range("F294").Value2 = 40664
$ws.Range("G294").Value = "Robert Parker"

$ws.Range("A295").Value = "Château Pavie"
$ws.Range("B295").Value = "Other"
$ws.Range("C295").Value = 2011
$ws.Range("D295").Value = 93
$ws.Range("E295").Value = 95
$ws.Range("F295").Value2 = 41000
$ws.Range("G295").Value = "Robert Parker"

$ws.Range("A296").Value = "Château Pavie"
$ws.Range("B296").Value = "Other"
$ws.Range("C296").Value = 2012
$ws.Range("D296").Value = 94
$ws.Range("E296").Value = 96
$ws.Range("F296").Value2 = 41365
$ws.Range("G296").Value = "Robert Parker"

$ws.Range("A297").Value = "Château Pavie"
$ws.Range("B297").Value = "Other"
$ws.Range("C297").Value = 2013
$ws.Range("D297").Value = 92
$ws.Range("E297").Value = 94
$ws.Range("F297").Value2 = 41852
$ws.Range("G297").Value = "Robert Parker"

$ws.Range("A298").Value = "Château Pavie"
$ws.Range("B298").Value = "Other"
$ws.Range("C298").Value = 2014
$ws.Range("D298").Value = 94
$ws.Range("E298").Value = 96
$ws.Range("F298").Value2 = 42095
$ws.Range("G298").Value = "Neal Martin"

$ws.Range("A299").Value = "Château Pavie"
$ws.Range("B299").Value = "Other"
$ws.Range("C299").Value = 2015
$ws.Range("D299").Value = 96
$ws.Range("E299").Value = 98
$ws.Range("F299").Value2 = 42461
$ws.Range("G299").Value = "Neal Martin"

$ws.Range("A300").Value = "Château Pavie"
$ws.Range("B300").Value = "Other"
$ws.Range("C300").Value = 2016
$ws.Range("D300").Value = 98
$ws.Range("E300").Value = 100
$ws.Range("F300").Value2 = 42826
$ws.Range("G300").Value = "Neal Martin"

$ws.Range("A301").Value = "Château Pavie"
$ws.Range("B301").Value = "Other"
$ws.Range("C301").Value = 2017
$ws.Range("D301").Value = 97
$ws.Range("E301").Value = 99
$ws.Range("F301").Value2 = 43191
$ws.Range("G301").Value = "Lisa Perrotti-Brown"

$ws.Range("A302").Value = "Château Pavie"
$ws.Range("B302").Value = "Other"
$ws.Range("C302").Value = 2018
$ws.Range("D302").Value = 97
$ws.Range("E302").Value = 100
$ws.Range("F302").Value2 = 43556
$ws.Range("G302").Value = "Lisa Perrotti-Brown"

# Update the view: scroll position and active selection moved to reflect the newly entered rows
[void]$excel.Goto($ws.Range("A261"))
[void]$ws.Range("A278").Select()
